# Refresh the cryptos snapshot (GitHub Actions scheduled update).
# coinranking.com price/volume(1h) figures change every run; row 40/41
# (Aave / MXToken) also swap rank order this time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells that look like plain numbers get a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# inline-string cells) instead of silently reformatting/rounding them
# (e.g. "0.1040" -> 0.104, "0.06890" -> 0.0689).

# Row 2: Bitcoin
$ws.Range("D2").Value = "30.855.68"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.936.82"
$ws.Range("E3").Value = "  -1.01%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5: BNB
$ws.Range("D5").Value = "'243.48"
$ws.Range("E5").Value = "  -1.22%  "

# Row 6: USDC
$ws.Range("E6").Value = "  -0.09%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.4898"
$ws.Range("E7").Value = "  -0.44%  "

# Row 8: Cardano
$ws.Range("D8").Value = "'0.2954"
$ws.Range("E8").Value = "  -0.78%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.06890"
$ws.Range("E9").Value = "  +0.61%  "

# Row 10: Solana
$ws.Range("E10").Value = "  +0.42%  "

# Row 11: Litecoin
$ws.Range("D11").Value = "'104.88"
$ws.Range("E11").Value = "  -2.83%  "

# Row 12: TRON
$ws.Range("D12").Value = "'0.07791"
$ws.Range("E12").Value = "  +0.32%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.935.37"
$ws.Range("E13").Value = "  -0.93%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'5.348"
$ws.Range("E14").Value = "  -2.22%  "

# Row 15: Polygon
$ws.Range("D15").Value = "'0.7011"
$ws.Range("E15").Value = "  -1.13%  "

# Row 16: BitcoinCash
$ws.Range("D16").Value = "'273.71"
$ws.Range("E16").Value = "  -3.11%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "30.846.77"
$ws.Range("E17").Value = "  -0.82%  "

# Row 18: ShibaInu
$ws.Range("D18").Value = "'0.000007716"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19: Avalanche
$ws.Range("D19").Value = "'13.06"
$ws.Range("E19").Value = "  -1.98%  "

# Row 20: Dai
$ws.Range("E20").Value = "  +0.00%  "

# Row 21: Uniswap
$ws.Range("D21").Value = "'5.575"
$ws.Range("E21").Value = "  +1.20%  "

# Row 22: WrappedliquidstakedEther2.0
$ws.Range("D22").Value = "2.188.19"
$ws.Range("E22").Value = "  -0.11%  "

# Row 23: BinanceUSD
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -0.17%  "

# Row 24: Chainlink
$ws.Range("D24").Value = "'6.535"
$ws.Range("E24").Value = "  +0.44%  "

# Row 25: Cosmos
$ws.Range("D25").Value = "'9.851"
$ws.Range("E25").Value = "  +0.10%  "

# Row 26: Monero
$ws.Range("D26").Value = "'166.06"
$ws.Range("E26").Value = "  -2.22%  "

# Row 27: EthereumClassic
$ws.Range("D27").Value = "'19.57"
$ws.Range("E27").Value = "  -2.53%  "

# Row 28: LidoDAOToken
$ws.Range("D28").Value = "'2.156"
$ws.Range("E28").Value = "  -2.55%  "

# Row 29: Stellar
$ws.Range("D29").Value = "'0.1040"
$ws.Range("E29").Value = "  -1.61%  "

# Row 30: Toncoin
$ws.Range("D30").Value = "'1.391"
$ws.Range("E30").Value = "  -2.78%  "

# Row 31: PancakeSwap
$ws.Range("D31").Value = "'1.560"
$ws.Range("E31").Value = "  -1.64%  "

# Row 32: Filecoin
$ws.Range("D32").Value = "'4.558"
$ws.Range("E32").Value = "  -0.59%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = "'4.373"
$ws.Range("E33").Value = "  -1.93%  "

# Row 34: Hedera
$ws.Range("D34").Value = "'0.04886"
$ws.Range("E34").Value = "  -1.87%  "

# Row 35: ImmutableX
$ws.Range("D35").Value = "'0.7614"
$ws.Range("E35").Value = "  +0.21%  "

# Row 36: ARBITRUM
$ws.Range("D36").Value = "'1.149"
$ws.Range("E36").Value = "  -3.11%  "

# Row 37: Frax
$ws.Range("E37").Value = "  -0.03%  "

# Row 38: HuobiToken
$ws.Range("D38").Value = "'2.732"
$ws.Range("E38").Value = "  -0.13%  "

# Row 39: VeChain
$ws.Range("D39").Value = "'0.02008"
$ws.Range("E39").Value = "  -1.51%  "

# Row 40: Aave
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'78.86"
$ws.Range("E40").Value = "  +5.74%  "

# Row 41: MXToken
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.658"
$ws.Range("E41").Value = "  -1.81%  "

# Row 42: FraxShare
$ws.Range("D42").Value = "'6.504"
$ws.Range("E42").Value = "  +0.02%  "

# Row 43: RenderToken
$ws.Range("D43").Value = "'2.086"
$ws.Range("E43").Value = "  -4.27%  "

# Row 44: TrustWalletToken
$ws.Range("D44").Value = "'0.9051"
$ws.Range("E44").Value = "  +2.27%  "

# Row 45: TheSandbox
$ws.Range("D45").Value = "'0.4436"
$ws.Range("E45").Value = "  -1.68%  "

# Row 46: Quant
$ws.Range("D46").Value = "'107.77"
$ws.Range("E46").Value = "  -1.52%  "

# Row 47: PaxDollar
$ws.Range("E47").Value = "  -0.08%  "

# Row 48: Aptos
$ws.Range("D48").Value = "'7.719"
$ws.Range("E48").Value = "  -5.92%  "

# Row 49: Maker
$ws.Range("D49").Value = "1.000.90"
$ws.Range("E49").Value = "  +2.28%  "

# Row 50: Algorand
$ws.Range("D50").Value = "'0.1250"
$ws.Range("E50").Value = "  -1.50%  "

# Row 51: Elrond
$ws.Range("D51").Value = "'36.21"
$ws.Range("E51").Value = "  +1.10%  "
